# The workbook tracks reviews; one row (the review + recovery-email pair
# for shmulmaor2@gmail.com) was removed from the sheet. Delete that
# worksheet row, which shifts everything below it up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Delete()
